$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 10.83367066666667
$ws.Range("N2").Value = 32.501012
$ws.Range("O2").Value = 0.1945663943642109
$ws.Range("P2").Value = 0.1945663943642109
$ws.Range("Q2").Value = 1.732307550823556
$ws.Range("R2").Value = 15.590767957412
$ws.Range("S2").Value = 0.005151674401392354
$ws.Range("T2").Value = 0.005151674401392356

$ws.Range("O3").Value = 0.5006586046969607
$ws.Range("P3").Value = 0.5006586046969607
$ws.Range("R3").Value = 40.118192852464
$ws.Range("S3").Value = 0.01325629806772314
$ws.Range("T3").Value = 0.01325629806772314

$ws.Range("M4").Value = 0.4511806666666667
$ws.Range("N4").Value = 1.353542
$ws.Range("O4").Value = 0.008102941119511068
$ws.Range("P4").Value = 0.00810294111951107
$ws.Range("Q4").Value = 0.07214393899355556
$ws.Range("R4").Value = 0.649295450942
$ws.Range("S4").Value = 0.0002145474015581241
$ws.Range("T4").Value = 0.0002145474015581242

$ws.Range("M5").Value = 16.146255
$ws.Range("N5").Value = 48.438765
$ws.Range("O5").Value = 0.2899773045068669
$ws.Range("P5").Value = 0.2899773045068669
$ws.Range("Q5").Value = 2.581791556585
$ws.Range("R5").Value = 23.236124009265
$ws.Range("S5").Value = 0.007677937711156808
$ws.Range("T5").Value = 0.007677937711156809

$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3727713333333333
$ws.Range("N6").Value = 1.118314
$ws.Range("O6").Value = 0.006694755312450519
$ws.Range("P6").Value = 0.006694755312450519
$ws.Range("Q6").Value = 0.05960626045711111
$ws.Range("R6").Value = 0.5364563441139999
$ws.Range("S6").Value = 0.0001772618528468802
$ws.Range("T6").Value = 0.0001772618528468803

$ws.Range("M7").Value = 10.83367066666667
$ws.Range("N7").Value = 32.501012
$ws.Range("O7").Value = 0.1945663943642109
$ws.Range("P7").Value = 0.1945663943642109
$ws.Range("Q7").Value = 63.69279656727468
$ws.Range("R7").Value = 573.2351691054721
$ws.Range("S7").Value = 0.1894147199628185
$ws.Range("T7").Value = 0.1894147199628186

$ws.Range("O8").Value = 0.5006586046969607
$ws.Range("P8").Value = 0.5006586046969607
$ws.Range("S8").Value = 0.4874023066292375
$ws.Range("T8").Value = 0.4874023066292376

$ws.Range("M9").Value = 0.4511806666666667
$ws.Range("N9").Value = 1.353542
$ws.Range("O9").Value = 0.008102941119511068
$ws.Range("P9").Value = 0.00810294111951107
$ws.Range("Q9").Value = 2.652559718794667
$ws.Range("R9").Value = 23.873037469152
$ws.Range("S9").Value = 0.007888393717952944
$ws.Range("T9").Value = 0.007888393717952946

$ws.Range("M10").Value = 16.146255
$ws.Range("N10").Value = 48.438765
$ws.Range("O10").Value = 0.2899773045068669
$ws.Range("P10").Value = 0.2899773045068669
$ws.Range("Q10").Value = 94.92628737576001
$ws.Range("R10").Value = 854.3365863818401
$ws.Range("S10").Value = 0.28229936679571
$ws.Range("T10").Value = 0.2822993667957101

$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3727713333333333
$ws.Range("N11").Value = 1.118314
$ws.Range("O11").Value = 0.006694755312450519
$ws.Range("P11").Value = 0.006694755312450519
$ws.Range("Q11").Value = 2.191579329909333
$ws.Range("R11").Value = 19.724213969184
$ws.Range("S11").Value = 0.006517493459603638
$ws.Range("T11").Value = 0.006517493459603639

